# "Correções Junto Com o professor"
# Adds a third column ("Area") with the average-area-per-stratum figures
# next to the existing "Estratos"/"Variavel" columns, and gives that new
# column's data rows (C3:C23) a centered + underlined look to set them
# apart from the header/first data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header for column C
$ws.Range("C1").Value = "Area"

# Area values per row (2..23)
$areaValues = @(
    14.4, # C2
    14.4, # C3
    14.4, # C4
    14.4, # C5
    14.4, # C6
    14.4, # C7
    14.4, # C8
    16.4, # C9
    16.4, # C10
    16.4, # C11
    16.4, # C12
    16.4, # C13
    16.4, # C14
    16.4, # C15
    16.4, # C16
    14.2, # C17
    14.2, # C18
    14.2, # C19
    14.2, # C20
    14.2, # C21
    14.2, # C22
    14.2  # C23
)

for ($i = 0; $i -lt $areaValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $areaValues[$i]
}

# Center-align the header and the whole new column's values
$ws.Range("C1:C23").HorizontalAlignment = -4108

# Underline the repeated-area rows (everything except the first data row)
# which gives them their own (3rd) cell style, distinct from the header
# and from C2.
$ws.Range("C3:C23").Font.Underline = 2

# Page setup tweaks that came along with the edit
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move/restore the active selection as recorded after the edit
$ws.Range("F14").Select()
